$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.527940154075623
$ws.Range("B1").Value = 3.667423009872437
$ws.Range("C1").Value = 5.82089900970459
$ws.Range("D1").Value = 1.422928214073181
$ws.Range("E1").Value = 0.8312093019485474
